# Refresh the scraped crypto price/volume figures (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column holds numeric-looking text (e.g. "8.200", "1.849.63"); keep it
# formatted as Text up front so COM does not silently coerce it to a Double
# (which would both reparse multi-dot values and drop trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.773.18"
$ws.Range("E2").Value = "  +1.77%  "

$ws.Range("D3").Value = "1.662.39"
$ws.Range("E3").Value = "  +1.87%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "329.85"
$ws.Range("E5").Value = "  +8.38%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").Value = "0.3641"
$ws.Range("E7").Value = "  +1.43%  "

$ws.Range("D8").Value = "47.28"
$ws.Range("E8").Value = "  +0.95%  "

$ws.Range("D9").Value = "0.3247"
$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("D10").Value = "1.135"
$ws.Range("E10").Value = "  +3.15%  "

$ws.Range("D11").Value = "0.07054"
$ws.Range("E11").Value = "  +2.85%  "

$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").Value = "6.059"
$ws.Range("E13").Value = "  +2.70%  "

$ws.Range("D14").Value = "19.47"
$ws.Range("E14").Value = "  +2.38%  "

$ws.Range("D15").Value = "1.666.59"
$ws.Range("E15").Value = "  +1.88%  "

$ws.Range("D16").Value = "6.577"
$ws.Range("E16").Value = "  +1.15%  "

$ws.Range("D17").Value = "0.00001046"
$ws.Range("E17").Value = "  +0.69%  "

$ws.Range("D18").Value = "0.06646"
$ws.Range("E18").Value = "  +2.24%  "

$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").Value = "78.38"
$ws.Range("E20").Value = "  +2.86%  "

$ws.Range("D21").Value = "5.914"
$ws.Range("E21").Value = "  +0.84%  "

$ws.Range("D22").Value = "15.75"
$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("D23").Value = "12.51"
$ws.Range("E23").Value = "  +4.82%  "

$ws.Range("D24").Value = "24.786.63"
$ws.Range("E24").Value = "  +2.03%  "

$ws.Range("D25").Value = "2.464"
$ws.Range("E25").Value = "  +2.94%  "

$ws.Range("D26").Value = "2.423"
$ws.Range("E26").Value = "  +5.19%  "

$ws.Range("D27").Value = "148.80"
$ws.Range("E27").Value = "  +3.53%  "

$ws.Range("D28").Value = "18.63"
$ws.Range("E28").Value = "  +0.70%  "

$ws.Range("D29").Value = "1.849.63"
$ws.Range("E29").Value = "  +1.82%  "

$ws.Range("D30").Value = "125.68"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("D31").Value = "1.163"
$ws.Range("E31").Value = "  +4.80%  "

$ws.Range("D32").Value = "4.066"
$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("D33").Value = "5.689"
$ws.Range("E33").Value = "  +1.21%  "

$ws.Range("D34").Value = "0.08485"
$ws.Range("E34").Value = "  +1.70%  "

$ws.Range("D35").Value = "1.641"
$ws.Range("E35").Value = "  -1.41%  "

$ws.Range("D36").Value = "12.12"
$ws.Range("E36").Value = "  -1.06%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "5.159"
$ws.Range("E37").Value = "  +1.43%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06173"
$ws.Range("E38").Value = "  +3.63%  "

$ws.Range("D39").Value = "0.02283"
$ws.Range("E39").Value = "  +3.86%  "

$ws.Range("D40").Value = "1.242"
$ws.Range("E40").Value = "  +3.78%  "

$ws.Range("D41").Value = "0.2086"
$ws.Range("E41").Value = "  +3.16%  "

$ws.Range("D42").Value = "8.200"
$ws.Range("E42").Value = "  +1.22%  "

$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("D44").Value = "0.5912"
$ws.Range("E44").Value = "  +1.87%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.41"
$ws.Range("E45").Value = "  +7.77%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.848"
$ws.Range("E46").Value = "  +3.91%  "

$ws.Range("D47").Value = "0.5663"
$ws.Range("E47").Value = "  +3.04%  "

$ws.Range("D48").Value = "125.54"
$ws.Range("E48").Value = "  +4.03%  "

$ws.Range("E49").Value = "  +1.85%  "

$ws.Range("D50").Value = "0.06964"
$ws.Range("E50").Value = "  +1.36%  "

$ws.Range("D51").Value = "1.192"
$ws.Range("E51").Value = "  +5.07%  "
